$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.232.88"
Set-TextValue $ws.Range("E2") "  +0.00%  "
Set-TextValue $ws.Range("D3") "1.901.99"
Set-TextValue $ws.Range("E3") "  -0.19%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "326.02"
Set-TextValue $ws.Range("E5") "  -0.54%  "
Set-TextValue $ws.Range("E6") "  -0.09%  "
Set-TextValue $ws.Range("D7") "0.4646"
Set-TextValue $ws.Range("E7") "  +0.44%  "
Set-TextValue $ws.Range("D8") "0.3917"
Set-TextValue $ws.Range("E8") "  -0.97%  "
Set-TextValue $ws.Range("D9") "0.07885"
Set-TextValue $ws.Range("E9") "  -1.02%  "
Set-TextValue $ws.Range("D10") "0.9901"
Set-TextValue $ws.Range("E10") "  -1.04%  "
Set-TextValue $ws.Range("D11") "21.79"
Set-TextValue $ws.Range("E11") "  -2.28%  "
Set-TextValue $ws.Range("D12") "1.945.59"
Set-TextValue $ws.Range("E12") "  +1.17%  "
Set-TextValue $ws.Range("D13") "7.075"
Set-TextValue $ws.Range("E13") "  -0.83%  "
Set-TextValue $ws.Range("D14") "5.742"
Set-TextValue $ws.Range("E14") "  -0.32%  "
Set-TextValue $ws.Range("D15") "0.06989"
Set-TextValue $ws.Range("E15") "  +0.53%  "
Set-TextValue $ws.Range("D16") "88.23"
Set-TextValue $ws.Range("E16") "  -0.60%  "
Set-TextValue $ws.Range("D17") "1.004"
Set-TextValue $ws.Range("E17") "  -0.02%  "
Set-TextValue $ws.Range("D18") "0.000009979"
Set-TextValue $ws.Range("E18") "  -1.26%  "
Set-TextValue $ws.Range("D19") "17.10"
Set-TextValue $ws.Range("E19") "  -0.43%  "
Set-TextValue $ws.Range("D21") "29.242.09"
Set-TextValue $ws.Range("E21") "  -0.08%  "
Set-TextValue $ws.Range("D22") "5.301"
Set-TextValue $ws.Range("E22") "  -1.14%  "
Set-TextValue $ws.Range("E23") "  -0.13%  "
Set-TextValue $ws.Range("D24") "2.126.15"
Set-TextValue $ws.Range("E24") "  -1.03%  "
Set-TextValue $ws.Range("D25") "2.110"
Set-TextValue $ws.Range("E25") "  +3.09%  "
Set-TextValue $ws.Range("D26") "156.22"
Set-TextValue $ws.Range("E26") "  -0.37%  "
Set-TextValue $ws.Range("D27") "19.41"
Set-TextValue $ws.Range("E27") "  -0.58%  "
Set-TextValue $ws.Range("D28") "5.972"
Set-TextValue $ws.Range("E28") "  +0.91%  "
Set-TextValue $ws.Range("D29") "118.61"
Set-TextValue $ws.Range("E29") "  -0.49%  "
Set-TextValue $ws.Range("D30") "1.884"
Set-TextValue $ws.Range("E30") "  -5.76%  "
Set-TextValue $ws.Range("D31") "0.09327"
Set-TextValue $ws.Range("E31") "  -0.89%  "
Set-TextValue $ws.Range("D32") "0.9021"
Set-TextValue $ws.Range("E32") "  -2.45%  "
Set-TextValue $ws.Range("D33") "5.257"
Set-TextValue $ws.Range("E33") "  -1.71%  "
Set-TextValue $ws.Range("D34") "1.325"
Set-TextValue $ws.Range("E34") "  -1.88%  "
Set-TextValue $ws.Range("D35") "3.191"
Set-TextValue $ws.Range("E35") "  -2.09%  "
Set-TextValue $ws.Range("E36") "  +1.11%  "
Set-TextValue $ws.Range("D37") "0.05772"
Set-TextValue $ws.Range("E37") "  -1.02%  "
Set-TextValue $ws.Range("D38") "0.02087"
Set-TextValue $ws.Range("E38") "  -1.05%  "
Set-TextValue $ws.Range("E39") "  -0.15%  "
Set-TextValue $ws.Range("D40") "7.722"
Set-TextValue $ws.Range("E40") "  -3.45%  "
Set-TextValue $ws.Range("D41") "0.5704"
Set-TextValue $ws.Range("E41") "  -0.95%  "
Set-TextValue $ws.Range("D42") "0.1790"
Set-TextValue $ws.Range("E42") "  -1.02%  "
Set-TextValue $ws.Range("D43") "9.713"
Set-TextValue $ws.Range("E43") "  -2.74%  "
Set-TextValue $ws.Range("D45") "0.5358"
Set-TextValue $ws.Range("E45") "  -1.35%  "
Set-TextValue $ws.Range("D46") "2.172"
Set-TextValue $ws.Range("E46") "  -2.44%  "
Set-TextValue $ws.Range("D47") "0.07017"
Set-TextValue $ws.Range("E47") "  -1.22%  "
Set-TextValue $ws.Range("D48") "1.853"
Set-TextValue $ws.Range("E48") "  -1.35%  "
Set-TextValue $ws.Range("D49") "2.579"
Set-TextValue $ws.Range("E49") "  +0.27%  "
Set-TextValue $ws.Range("D50") "113.26"
Set-TextValue $ws.Range("E50") "  +1.10%  "
Set-TextValue $ws.Range("D51") "1.062"
Set-TextValue $ws.Range("E51") "  +0.21%  "
